# Update existing rows 2-13 and add new rows 14-17 for the Fgf1-Fgfr2 sheet
# with the refreshed TPM-based NATMI values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7144740000000001
$ws.Range("N2").Value = 2.143422
$ws.Range("O2").Value = 0.138796410342318
$ws.Range("P2").Value = 0.138796410342318
$ws.Range("Q2").Value = 1.288262591766
$ws.Range("R2").Value = 11.594363325894
$ws.Range("S2").Value = 0.02421839861190604
$ws.Range("T2").Value = 0.02421839861190604

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.140873
$ws.Range("N3").Value = 12.422619
$ws.Range("O3").Value = 0.8044215857867821
$ws.Range("P3").Value = 0.8044215857867821
$ws.Range("Q3").Value = 7.466376359606999
$ws.Range("R3").Value = 67.197387236463
$ws.Range("S3").Value = 0.1403624385425911
$ws.Range("T3").Value = 0.140362438542591

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2847646666666667
$ws.Range("N4").Value = 0.8542940000000001
$ws.Range("O4").Value = 0.05531945672713084
$ws.Range("P4").Value = 0.05531945672713083
$ws.Range("Q4").Value = 0.513456987270889
$ws.Range("R4").Value = 4.621112885438
$ws.Range("S4").Value = 0.009652617461125092
$ws.Range("T4").Value = 0.009652617461125087

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.803092333333333
$ws.Range("H5").Value = 5.409276999999999
$ws.Range("I5").Value = 0.1744886524959502
$ws.Range("J5").Value = 0.1744886524959502
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.007528666666666667
$ws.Range("N5").Value = 0.022586
$ws.Range("O5").Value = 0.00146254714376898
$ws.Range("P5").Value = 0.00146254714376898
$ws.Range("Q5").Value = 0.01357488114688889
$ws.Range("R5").Value = 0.122173930322
$ws.Range("S5").Value = 0.0002551978803280501
$ws.Range("T5").Value = 0.0002551978803280501

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.531154666666667
$ws.Range("H6").Value = 13.593464
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7144740000000001
$ws.Range("N6").Value = 2.143422
$ws.Range("O6").Value = 0.138796410342318
$ws.Range("P6").Value = 0.138796410342318
$ws.Range("Q6").Value = 3.237392199312
$ws.Range("R6").Value = 29.136529793808
$ws.Range("S6").Value = 0.06086061587687131
$ws.Range("T6").Value = 0.0608606158768713

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.531154666666667
$ws.Range("H7").Value = 13.593464
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.140873
$ws.Range("N7").Value = 12.422619
$ws.Range("O7").Value = 0.8044215857867821
$ws.Range("P7").Value = 0.8044215857867821
$ws.Range("Q7").Value = 18.762936018024
$ws.Range("R7").Value = 168.866424162216
$ws.Range("S7").Value = 0.3527295339619184
$ws.Range("T7").Value = 0.3527295339619184

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.531154666666667
$ws.Range("H8").Value = 13.593464
$ws.Range("I8").Value = 0.4384883998568034
$ws.Range("J8").Value = 0.4384883998568034
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2847646666666667
$ws.Range("N8").Value = 0.8542940000000001
$ws.Range("O8").Value = 0.05531945672713084
$ws.Range("P8").Value = 0.05531945672713083
$ws.Range("Q8").Value = 1.290312748268445
$ws.Range("R8").Value = 11.612814734416
$ws.Range("S8").Value = 0.02425694006122728
$ws.Range("T8").Value = 0.02425694006122727

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.531154666666667
$ws.Range("H9").Value = 13.593464
$ws.Range("I9").Value = 0.4384883998568034
$ws.Range("J9").Value = 0.4384883998568034
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.007528666666666667
$ws.Range("N9").Value = 0.022586
$ws.Range("O9").Value = 0.00146254714376898
$ws.Range("P9").Value = 0.00146254714376898
$ws.Range("Q9").Value = 0.03411355310044445
$ws.Range("R9").Value = 0.307021977904
$ws.Range("S9").Value = 0.0006413099567863983
$ws.Range("T9").Value = 0.0006413099567863981

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7144740000000001
$ws.Range("N10").Value = 2.143422
$ws.Range("O10").Value = 0.138796410342318
$ws.Range("P10").Value = 0.138796410342318
$ws.Range("Q10").Value = 2.783308963086001
$ws.Range("R10").Value = 25.04978066777401
$ws.Range("S10").Value = 0.05232418169940271
$ws.Range("T10").Value = 0.05232418169940271

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgf1"
$ws.Range("C11").Value = "Fgfr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.895605666666667
$ws.Range("H11").Value = 11.686817
$ws.Range("I11").Value = 0.3769851220961256
$ws.Range("J11").Value = 0.3769851220961256
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.140873
$ws.Range("N11").Value = 12.422619
$ws.Range("O11").Value = 0.8044215857867821
$ws.Range("P11").Value = 0.8044215857867821
$ws.Range("Q11").Value = 16.131208323747
$ws.Range("R11").Value = 145.180874913723
$ws.Range("S11").Value = 0.303254969734589
$ws.Range("T11").Value = 0.303254969734589

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Fgf1"
$ws.Range("C12").Value = "Fgfr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.895605666666667
$ws.Range("H12").Value = 11.686817
$ws.Range("I12").Value = 0.3769851220961256
$ws.Range("J12").Value = 0.3769851220961256
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2847646666666667
$ws.Range("N12").Value = 0.8542940000000001
$ws.Range("O12").Value = 0.05531945672713084
$ws.Range("P12").Value = 0.05531945672713083
$ws.Range("Q12").Value = 1.109330849133111
$ws.Range("R12").Value = 9.983977642198003
$ws.Range("S12").Value = 0.02085461214856876
$ws.Range("T12").Value = 0.02085461214856875

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Fgf1"
$ws.Range("C13").Value = "Fgfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.895605666666667
$ws.Range("H13").Value = 11.686817
$ws.Range("I13").Value = 0.3769851220961256
$ws.Range("J13").Value = 0.3769851220961256
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.007528666666666667
$ws.Range("N13").Value = 0.022586
$ws.Range("O13").Value = 0.00146254714376898
$ws.Range("P13").Value = 0.00146254714376898
$ws.Range("Q13").Value = 0.02932871652911111
$ws.Range("R13").Value = 0.263958448762
$ws.Range("S13").Value = 0.0005513585135650887
$ws.Range("T13").Value = 0.0005513585135650886

# Row 14: Resolving-Mac -> ECs
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Fgf1"
$ws.Range("C14").Value = "Fgfr2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1037266666666667
$ws.Range("H14").Value = 0.31118
$ws.Range("I14").Value = 0.01003782555112075
$ws.Range("J14").Value = 0.01003782555112075
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.7144740000000001
$ws.Range("N14").Value = 2.143422
$ws.Range("O14").Value = 0.138796410342318
$ws.Range("P14").Value = 0.138796410342318
$ws.Range("Q14").Value = 0.07411000644000001
$ws.Range("R14").Value = 0.6669900579600001
$ws.Range("S14").Value = 0.00139321415413796
$ws.Range("T14").Value = 0.00139321415413796

# Row 15: Resolving-Mac -> FAPs
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Fgf1"
$ws.Range("C15").Value = "Fgfr2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1037266666666667
$ws.Range("H15").Value = 0.31118
$ws.Range("I15").Value = 0.01003782555112075
$ws.Range("J15").Value = 0.01003782555112075
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.140873
$ws.Range("N15").Value = 12.422619
$ws.Range("O15").Value = 0.8044215857867821
$ws.Range("P15").Value = 0.8044215857867821
$ws.Range("Q15").Value = 0.42951895338
$ws.Range("R15").Value = 3.865670580420001
$ws.Range("S15").Value = 0.008074643547683636
$ws.Range("T15").Value = 0.008074643547683634

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Fgf1"
$ws.Range("C16").Value = "Fgfr2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1037266666666667
$ws.Range("H16").Value = 0.31118
$ws.Range("I16").Value = 0.01003782555112075
$ws.Range("J16").Value = 0.01003782555112075
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2847646666666667
$ws.Range("N16").Value = 0.8542940000000001
$ws.Range("O16").Value = 0.05531945672713084
$ws.Range("P16").Value = 0.05531945672713083
$ws.Range("Q16").Value = 0.02953768965777779
$ws.Range("R16").Value = 0.2658392069200001
$ws.Range("S16").Value = 0.0005552870562097128
$ws.Range("T16").Value = 0.0005552870562097125

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Fgf1"
$ws.Range("C17").Value = "Fgfr2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.1037266666666667
$ws.Range("H17").Value = 0.31118
$ws.Range("I17").Value = 0.01003782555112075
$ws.Range("J17").Value = 0.01003782555112075
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.007528666666666667
$ws.Range("N17").Value = 0.022586
$ws.Range("O17").Value = 0.00146254714376898
$ws.Range("P17").Value = 0.00146254714376898
$ws.Range("Q17").Value = 0.0007809234977777779
$ws.Range("R17").Value = 0.00702831148
$ws.Range("S17").Value = 0.00001468079308944294
$ws.Range("T17").Value = 0.00001468079308944294

